$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 13326.68
$ws.Range("I33").Value = 14871.272
$ws.Range("K33").Value = 14871.272
$ws.Range("M33").Value = -14642.272

$ws.Range("H43").Value = 2415.3333
$ws.Range("I43").Value = 1623.5
$ws.Range("K43").Value = 1623.5
$ws.Range("M43").Value = -1554.5

$ws.Range("H116").Value = 3942.8125
$ws.Range("I116").Value = 3214.8333
$ws.Range("J116").Value = 4379.6
$ws.Range("K116").Value = 3214.8333
$ws.Range("L116").Value = 4379.6
$ws.Range("M116").Value = 227.1667000000002
$ws.Range("N116").Value = -11263.6

$ws.Range("H129").Value = 1304174.9
$ws.Range("I129").Value = 1598.3
$ws.Range("K129").Value = 4794.9
$ws.Range("M129").Value = 205.1000000000004

$ws.Range("H132").Value = 4138.3335
$ws.Range("I132").Value = 4257.5806
$ws.Range("J132").Value = 3399
$ws.Range("K132").Value = 12772.7418
$ws.Range("L132").Value = 10197
$ws.Range("M132").Value = -10242.7418
$ws.Range("N132").Value = -15257

$ws.Range("H137").Value = 2586.7856
$ws.Range("J137").Value = 2914.8
$ws.Range("L137").Value = 8744.400000000001
$ws.Range("N137").Value = -13844.4

$ws.Range("H138").Value = 5869.7144
$ws.Range("J138").Value = 5813.1577
$ws.Range("L138").Value = 17439.4731
$ws.Range("N138").Value = -27719.4731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 16262.6
$ws.Range("I50").Value = 5439.3335
$ws.Range("K50").Value = 5439.3335
$ws.Range("M50").Value = -4725.3335

$ws.Range("H102").Value = 5081.727
$ws.Range("I102").Value = 3612.375
$ws.Range("K102").Value = 3612.375
$ws.Range("M102").Value = -1990.375

$ws.Range("H132").Value = 4282.064
$ws.Range("I132").Value = 3459.65
$ws.Range("K132").Value = 10378.95
$ws.Range("M132").Value = -7848.950000000001

$ws.Range("H135").Value = 116428.29
$ws.Range("J135").Value = 116428.29
$ws.Range("L135").Value = 116428.29
$ws.Range("N135").Value = -126568.29

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 696.5454999999999
$ws.Range("I64").Value = 1018.3333
$ws.Range("K64").Value = 1018.3333
$ws.Range("M64").Value = -793.3333

$ws.Range("H67").Value = 696.5454999999999
$ws.Range("I67").Value = 1018.3333
$ws.Range("K67").Value = 1018.3333
$ws.Range("M67").Value = -238.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6886
$ws.Range("I134").Value = 3097.1538
$ws.Range("K134").Value = 9291.4614
$ws.Range("M134").Value = -6756.4614

$ws.Range("H141").Value = 48333.332
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 52500
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 52500
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -62860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11629838
$ws.Range("J68").Value = 1958.08
$ws.Range("L68").Value = 5874.24
$ws.Range("N68").Value = -7496.24

$ws.Range("H71").Value = 11629838
$ws.Range("J71").Value = 1958.08
$ws.Range("L71").Value = 17622.72
$ws.Range("N71").Value = -25734.72

$ws.Range("H103").Value = 499.66666
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 549.5
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 1648.5
$ws.Range("M103").Value = -321
$ws.Range("N103").Value = -3406.5

$ws.Range("H124").Value = 3000
$ws.Range("I124").Value = 1000
$ws.Range("J124").Value = 5000
$ws.Range("K124").Value = 3000
$ws.Range("L124").Value = 15000
$ws.Range("M124").Value = 1910
$ws.Range("N124").Value = -24820

$ws.Range("H132").Value = 62501040
$ws.Range("J132").Value = 1544.5
$ws.Range("L132").Value = 13900.5
$ws.Range("N132").Value = -18960.5

$ws.Range("H140").Value = 2048.8572
$ws.Range("I140").Value = 1026.2858
$ws.Range("K140").Value = 3078.8574
$ws.Range("M140").Value = 2101.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6056.2
$ws.Range("J29").Value = 11500
$ws.Range("L29").Value = 11500
$ws.Range("N29").Value = -12080

$ws.Range("H31").Value = 3002.25
$ws.Range("I31").Value = 3002.25
$ws.Range("K31").Value = 3002.25
$ws.Range("M31").Value = -2710.25

$ws.Range("H37").Value = 3002.25
$ws.Range("I37").Value = 3002.25
$ws.Range("K37").Value = 3002.25
$ws.Range("M37").Value = -2725.25

$ws.Range("H132").Value = 3666.96
$ws.Range("I132").Value = 3758.0454
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 11274.1362
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -8744.136200000001
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2650.9565
$ws.Range("I22").Value = 1859.4286
$ws.Range("K22").Value = 1859.4286
$ws.Range("M22").Value = -1564.4286

$ws.Range("H27").Value = 2650.9565
$ws.Range("I27").Value = 1859.4286
$ws.Range("K27").Value = 1859.4286
$ws.Range("M27").Value = -1752.4286

$ws.Range("H46").Value = 5525.6294
$ws.Range("J46").Value = 2986.75
$ws.Range("L46").Value = 2986.75
$ws.Range("N46").Value = -3362.75

$ws.Range("H93").Value = 12764.27
$ws.Range("I93").Value = 4939.2383
$ws.Range("K93").Value = 4939.2383
$ws.Range("M93").Value = -3691.2383

$ws.Range("H132").Value = 20429
$ws.Range("I132").Value = 23000.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 69001.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -66471.5
$ws.Range("N132").Value = -20060

$ws.Range("H136").Value = 6919.4
$ws.Range("I136").Value = 1532.3334
$ws.Range("K136").Value = 4597.0002
$ws.Range("M136").Value = -2047.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H32").Value = 4500
$ws.Range("I32").Value = 4500
$ws.Range("K32").Value = 4500
$ws.Range("M32").Value = -4183

$ws.Range("H55").Value = 1363.25
$ws.Range("I55").Value = 950
$ws.Range("J55").Value = 1776.5
$ws.Range("K55").Value = 950
$ws.Range("L55").Value = 1776.5
$ws.Range("M55").Value = -673
$ws.Range("N55").Value = -2330.5

$ws.Range("H70").Value = 37875.25
$ws.Range("I70").Value = 30001
$ws.Range("K70").Value = 30001
$ws.Range("M70").Value = -29686

$ws.Range("H73").Value = 37875.25
$ws.Range("I73").Value = 30001
$ws.Range("K73").Value = 30001
$ws.Range("M73").Value = -28909

$ws.Range("H136").Value = 6616.3228
$ws.Range("I136").Value = 4933.8096
$ws.Range("K136").Value = 14801.4288
$ws.Range("M136").Value = -12251.4288
